# Rebuild the logical-schema listing: each paragraph gets explicit
# "majorHAnsi" theme-font / 12pt (sz 24 half-points) run formatting, the
# table definitions are rewritten with their real column lists, and two
# new entities (Pedido, ItensPedidos) are appended.

$d = $word.ActiveDocument

# ---- helpers -------------------------------------------------------

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPrXml = '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

function Esc([string]$s) {
    return $s.Replace('&', '&amp;').Replace('<', '&lt;').Replace('>', '&gt;')
}

# Builds a <w:r> run with the standard paragraph formatting applied.
function Run([string]$text) {
    $needsSpace = ($text -ne $text.Trim()) -or ($text -eq '')
    $spaceAttr = ''
    if ($needsSpace) { $spaceAttr = ' xml:space="preserve"' }
    return '<w:r>' + $rPrXml + '<w:t' + $spaceAttr + '>' + (Esc $text) + '</w:t></w:r>'
}

function ProofStart { '<w:proofErr w:type="spellStart"/>' }
function ProofEnd { '<w:proofErr w:type="spellEnd"/>' }
function GramStart { '<w:proofErr w:type="gramStart"/>' }
function GramEnd { '<w:proofErr w:type="gramEnd"/>' }

# Builds a full <w:p> whose pPr also carries the standard run formatting,
# from an array of fragments (plain strings become runs automatically;
# strings starting with "<w:" are inserted verbatim, e.g. proofErr marks).
function Paragraph([string[]]$parts) {
    $inner = ''
    foreach ($part in $parts) {
        if ($part.StartsWith('<w:')) {
            $inner += $part
        } else {
            $inner += (Run $part)
        }
    }
    return '<w:p><w:pPr>' + $rPrXml + '</w:pPr>' + $inner + '</w:p>'
}

# ---- paragraph contents ---------------------------------------------

$p1 = Paragraph @(
    'Funcionario', ' ', '(',
    'codFuncionario, nome, cargo, sexo, email, dataNascimento',
    (GramStart), ', ', ')', (GramEnd)
)

$p2 = Paragraph @(
    'Compra', ' ', '(codCompra)', '     '
)

$p3 = Paragraph @(
    'Produto', ' ', '(', 'codProduto, nome, unidadeMedida, ', 'valor', ')'
)

$p4 = Paragraph @(
    'Loja', ' ', '(', 'codLoja, nome, endereço, cnpj', ')'
)

$p5 = Paragraph @(
    'Pedido ', '(',
    (ProofStart), 'cod', 'Pedido', (ProofEnd),
    ', ', 'codFunciona', 'rio, valorTotal, data) '
)

$p6 = Paragraph @(
    (ProofStart), 'ItensPedidos', (ProofEnd),
    ' (',
    (ProofStart), 'codPedido', (ProofEnd),
    ', codProduto, ',
    (ProofStart), 'preco', (ProofEnd),
    ', quantidade)',
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
)

$bodyXml = $p1 + $p2 + $p3 + $p4 + $p5 + $p6

# ---- apply via a WordprocessingML package fragment -------------------
# Range.InsertXML replaces the contents of the range it is called on, so
# targeting the whole document body rewrites every paragraph in one shot.

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'

$d.Content.InsertXML($packageXml)

Write-Host "Schema paragraphs rewritten:" $d.Paragraphs.Count
